$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 16.33975766666667
$ws.Range("H2").Value = 49.01927300000001
$ws.Range("I2").Value = 0.2979519994155143
$ws.Range("J2").Value = 0.2979519994155143
$ws.Range("M2").Value = 1.174933333333333
$ws.Range("N2").Value = 3.5248
$ws.Range("O2").Value = 0.01171850713626266
$ws.Range("P2").Value = 0.01171850713626266
$ws.Range("Q2").Value = 19.19812594115556
$ws.Range("R2").Value = 172.7831334704
$ws.Range("S2").Value = 0.003491552631414431
$ws.Range("T2").Value = 0.003491552631414431
$ws.Range("G3").Value = 16.33975766666667
$ws.Range("H3").Value = 49.01927300000001
$ws.Range("I3").Value = 0.2979519994155143
$ws.Range("J3").Value = 0.2979519994155143
$ws.Range("O3").Value = 0.2743256641287217
$ws.Range("P3").Value = 0.2743256641287218
$ws.Range("Q3").Value = 449.4206119939249
$ws.Range("R3").Value = 4044.785507945325
$ws.Range("S3").Value = 0.08173588011814147
$ws.Range("T3").Value = 0.08173588011814149
$ws.Range("G4").Value = 16.33975766666667
$ws.Range("H4").Value = 49.01927300000001
$ws.Range("I4").Value = 0.2979519994155143
$ws.Range("J4").Value = 0.2979519994155143
$ws.Range("M4").Value = 39.361408
$ws.Range("N4").Value = 118.084224
$ws.Range("O4").Value = 0.3925813724534833
$ws.Range("P4").Value = 0.3925813724534833
$ws.Range("Q4").Value = 643.1558681387947
$ws.Range("R4").Value = 5788.402813249153
$ws.Range("S4").Value = 0.1169704048558021
$ws.Range("T4").Value = 0.1169704048558021
$ws.Range("G5").Value = 16.33975766666667
$ws.Range("H5").Value = 49.01927300000001
$ws.Range("I5").Value = 0.2979519994155143
$ws.Range("J5").Value = 0.2979519994155143
$ws.Range("M5").Value = 32.221985
$ws.Range("N5").Value = 96.665955
$ws.Range("O5").Value = 0.3213744562815322
$ws.Range("P5").Value = 0.3213744562815322
$ws.Range("Q5").Value = 526.4994264389683
$ws.Range("R5").Value = 4738.494837950715
$ws.Range("S5").Value = 0.09575416181015631
$ws.Range("T5").Value = 0.09575416181015632
$ws.Range("I6").Value = 0.2656466977818992
$ws.Range("J6").Value = 0.2656466977818992
$ws.Range("M6").Value = 1.174933333333333
$ws.Range("N6").Value = 3.5248
$ws.Range("O6").Value = 0.01171850713626266
$ws.Range("P6").Value = 0.01171850713626266
$ws.Range("Q6").Value = 17.11657840817778
$ws.Range("R6").Value = 154.0492056736
$ws.Range("S6").Value = 0.003112982723681795
$ws.Range("T6").Value = 0.003112982723681794
$ws.Range("I7").Value = 0.2656466977818992
$ws.Range("J7").Value = 0.2656466977818992
$ws.Range("O7").Value = 0.2743256641287217
$ws.Range("P7").Value = 0.2743256641287218
$ws.Range("S7").Value = 0.07287370679262133
$ws.Range("T7").Value = 0.07287370679262134
$ws.Range("I8").Value = 0.2656466977818992
$ws.Range("J8").Value = 0.2656466977818992
$ws.Range("M8").Value = 39.361408
$ws.Range("N8").Value = 118.084224
$ws.Range("O8").Value = 0.3925813724534833
$ws.Range("P8").Value = 0.3925813724534833
$ws.Range("Q8").Value = 573.4220037632853
$ws.Range("R8").Value = 5160.798033869568
$ws.Range("S8").Value = 0.1042879452029537
$ws.Range("T8").Value = 0.1042879452029537
$ws.Range("I9").Value = 0.2656466977818992
$ws.Range("J9").Value = 0.2656466977818992
$ws.Range("M9").Value = 32.221985
$ws.Range("N9").Value = 96.665955
$ws.Range("O9").Value = 0.3213744562815322
$ws.Range("P9").Value = 0.3213744562815322
$ws.Range("Q9").Value = 469.4139804127565
$ws.Range("R9").Value = 4224.72582371481
$ws.Range("S9").Value = 0.08537206306264235
$ws.Range("T9").Value = 0.08537206306264236
$ws.Range("G10").Value = 22.15292366666667
$ws.Range("H10").Value = 66.458771
$ws.Range("I10").Value = 0.4039538427701242
$ws.Range("J10").Value = 0.4039538427701242
$ws.Range("M10").Value = 1.174933333333333
$ws.Range("N10").Value = 3.5248
$ws.Range("O10").Value = 0.01171850713626266
$ws.Range("P10").Value = 0.01171850713626266
$ws.Range("Q10").Value = 26.02820844675556
$ws.Range("R10").Value = 234.2538760208
$ws.Range("S10").Value = 0.004733735989222425
$ws.Range("T10").Value = 0.004733735989222424
$ws.Range("G11").Value = 22.15292366666667
$ws.Range("H11").Value = 66.458771
$ws.Range("I11").Value = 0.4039538427701242
$ws.Range("J11").Value = 0.4039538427701242
$ws.Range("O11").Value = 0.2743256641287217
$ws.Range("P11").Value = 0.2743256641287218
$ws.Range("Q11").Value = 609.3101693936609
$ws.Range("R11").Value = 5483.791524542949
$ws.Range("S11").Value = 0.1108149061952636
$ws.Range("T11").Value = 0.1108149061952636
$ws.Range("G12").Value = 22.15292366666667
$ws.Range("H12").Value = 66.458771
$ws.Range("I12").Value = 0.4039538427701242
$ws.Range("J12").Value = 0.4039538427701242
$ws.Range("M12").Value = 39.361408
$ws.Range("N12").Value = 118.084224
$ws.Range("O12").Value = 0.3925813724534833
$ws.Range("P12").Value = 0.3925813724534833
$ws.Range("Q12").Value = 871.9702668365228
$ws.Range("R12").Value = 7847.732401528704
$ws.Range("S12").Value = 0.158584754002554
$ws.Range("T12").Value = 0.158584754002554
$ws.Range("G13").Value = 22.15292366666667
$ws.Range("H13").Value = 66.458771
$ws.Range("I13").Value = 0.4039538427701242
$ws.Range("J13").Value = 0.4039538427701242
$ws.Range("M13").Value = 32.221985
$ws.Range("N13").Value = 96.665955
$ws.Range("O13").Value = 0.3213744562815322
$ws.Range("P13").Value = 0.3213744562815322
$ws.Range("Q13").Value = 713.8111740934783
$ws.Range("R13").Value = 6424.300566841304
$ws.Range("S13").Value = 0.1298204465830842
$ws.Range("T13").Value = 0.1298204465830842
$ws.Range("G14").Value = 1.779426333333333
$ws.Range("H14").Value = 5.338279
$ws.Range("I14").Value = 0.03244746003246218
$ws.Range("J14").Value = 0.03244746003246217
$ws.Range("M14").Value = 1.174933333333333
$ws.Range("N14").Value = 3.5248
$ws.Range("O14").Value = 0.01171850713626266
$ws.Range("P14").Value = 0.01171850713626266
$ws.Range("Q14").Value = 2.090707313244445
$ws.Range("R14").Value = 18.8163658192
$ws.Range("S14").Value = 0.0003802357919440054
$ws.Range("T14").Value = 0.0003802357919440052
$ws.Range("G15").Value = 1.779426333333333
$ws.Range("H15").Value = 5.338279
$ws.Range("I15").Value = 0.03244746003246218
$ws.Range("J15").Value = 0.03244746003246217
$ws.Range("O15").Value = 0.2743256641287217
$ws.Range("P15").Value = 0.2743256641287218
$ws.Range("Q15").Value = 48.94263966693911
$ws.Range("R15").Value = 440.483757002452
$ws.Range("S15").Value = 0.008901171022695342
$ws.Range("T15").Value = 0.008901171022695342
$ws.Range("G16").Value = 1.779426333333333
$ws.Range("H16").Value = 5.338279
$ws.Range("I16").Value = 0.03244746003246218
$ws.Range("J16").Value = 0.03244746003246217
$ws.Range("M16").Value = 39.361408
$ws.Range("N16").Value = 118.084224
$ws.Range("O16").Value = 0.3925813724534833
$ws.Range("P16").Value = 0.3925813724534833
$ws.Range("Q16").Value = 70.04072591227734
$ws.Range("R16").Value = 630.3665332104961
$ws.Range("S16").Value = 0.01273826839217355
$ws.Range("T16").Value = 0.01273826839217355
$ws.Range("G17").Value = 1.779426333333333
$ws.Range("H17").Value = 5.338279
$ws.Range("I17").Value = 0.03244746003246218
$ws.Range("J17").Value = 0.03244746003246217
$ws.Range("M17").Value = 32.221985
$ws.Range("N17").Value = 96.665955
$ws.Range("O17").Value = 0.3213744562815322
$ws.Range("P17").Value = 0.3213744562815322
$ws.Range("Q17").Value = 57.33664862127166
$ws.Range("R17").Value = 516.029837591445
$ws.Range("S17").Value = 0.1042778482564928
$ws.Range("T17").Value = 0.01042778482564928
